# Generate and save output file after processing
# Insert 3 new columns for general_college_subjects.history / .electives / .cs
# right before the existing general_college_subjects.arts column (R), shifting
# everything from R onward three columns to the right, then populate the new
# header cells + data cells and refresh the descriptive text in D2:J2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank columns at R:T (old R "arts" and everything after shifts to U:AH)
$ws.Range("R1:T1").EntireColumn.Insert()

# New header labels for the inserted columns
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data values for row 2 in the inserted columns
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Updated descriptive values (previously "Unknown" placeholders)
$ws.Range("D2").Value = "considered"
$ws.Range("E2").Value = "important"
$ws.Range("F2").Value = "not considered"
$ws.Range("G2").Value = "very important"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "important"
$ws.Range("J2").Value = "important"

$wb.Save()
